$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$cleanStyle = $ws.Range("B2").Style

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.845.60'
$ws.Range('D2').Style = $cleanStyle
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.87%  '
$ws.Range('E2').Style = $cleanStyle
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.209.93'
$ws.Range('D3').Style = $cleanStyle
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.73%  '
$ws.Range('E3').Style = $cleanStyle
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('E4').Style = $cleanStyle
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.16'
$ws.Range('D5').Style = $cleanStyle
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -2.82%  '
$ws.Range('E5').Style = $cleanStyle
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.26%  '
$ws.Range('E6').Style = $cleanStyle
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '73.13'
$ws.Range('D7').Style = $cleanStyle
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -1.88%  '
$ws.Range('E7').Style = $cleanStyle
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.17%  '
$ws.Range('E8').Style = $cleanStyle
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.12%  '
$ws.Range('E9').Style = $cleanStyle
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.51'
$ws.Range('D10').Style = $cleanStyle
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.71%  '
$ws.Range('E10').Style = $cleanStyle
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0950'
$ws.Range('D11').Style = $cleanStyle
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.42%  '
$ws.Range('E11').Style = $cleanStyle
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.06'
$ws.Range('D12').Style = $cleanStyle
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.00%  '
$ws.Range('E12').Style = $cleanStyle
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.61%  '
$ws.Range('E13').Style = $cleanStyle
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.545.93'
$ws.Range('D14').Style = $cleanStyle
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.41%  '
$ws.Range('E14').Style = $cleanStyle
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.25'
$ws.Range('D15').Style = $cleanStyle
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -2.31%  '
$ws.Range('E15').Style = $cleanStyle
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.835'
$ws.Range('D16').Style = $cleanStyle
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -2.35%  '
$ws.Range('E16').Style = $cleanStyle
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.202.98'
$ws.Range('D17').Style = $cleanStyle
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.00%  '
$ws.Range('E17').Style = $cleanStyle
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '41.808.96'
$ws.Range('D18').Style = $cleanStyle
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.75%  '
$ws.Range('E18').Style = $cleanStyle
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +6.94%  '
$ws.Range('E19').Style = $cleanStyle
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.99'
$ws.Range('D20').Style = $cleanStyle
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.70%  '
$ws.Range('E20').Style = $cleanStyle
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.14'
$ws.Range('D21').Style = $cleanStyle
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('E21').Style = $cleanStyle
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +22.12%  '
$ws.Range('E22').Style = $cleanStyle
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '229.07'
$ws.Range('D23').Style = $cleanStyle
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.45%  '
$ws.Range('E23').Style = $cleanStyle
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -6.65%  '
$ws.Range('E24').Style = $cleanStyle
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.76'
$ws.Range('D25').Style = $cleanStyle
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +3.67%  '
$ws.Range('E25').Style = $cleanStyle
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.20%  '
$ws.Range('E26').Style = $cleanStyle
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.27%  '
$ws.Range('E27').Style = $cleanStyle
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.33%  '
$ws.Range('E28').Style = $cleanStyle
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.19'
$ws.Range('D29').Style = $cleanStyle
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -2.14%  '
$ws.Range('E29').Style = $cleanStyle
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '167.10'
$ws.Range('D30').Style = $cleanStyle
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.31%  '
$ws.Range('E31').Style = $cleanStyle
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.55'
$ws.Range('D32').Style = $cleanStyle
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +6.02%  '
$ws.Range('E32').Style = $cleanStyle
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0793'
$ws.Range('D33').Style = $cleanStyle
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -4.31%  '
$ws.Range('E33').Style = $cleanStyle
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'InjectiveProtocol'
$ws.Range('B34').Style = $cleanStyle
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('C34').Style = $cleanStyle
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '29.60'
$ws.Range('D34').Style = $cleanStyle
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -4.83%  '
$ws.Range('E34').Style = $cleanStyle
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'Stellar'
$ws.Range('B35').Style = $cleanStyle
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('C35').Style = $cleanStyle
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.124'
$ws.Range('D35').Style = $cleanStyle
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.81%  '
$ws.Range('E35').Style = $cleanStyle
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -10.60%  '
$ws.Range('E36').Style = $cleanStyle
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -5.09%  '
$ws.Range('E37').Style = $cleanStyle
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0299'
$ws.Range('D38').Style = $cleanStyle
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -5.30%  '
$ws.Range('E38').Style = $cleanStyle
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '13.79'
$ws.Range('D39').Style = $cleanStyle
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.74%  '
$ws.Range('E39').Style = $cleanStyle
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '65.60'
$ws.Range('D40').Style = $cleanStyle
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +4.60%  '
$ws.Range('E40').Style = $cleanStyle
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.11'
$ws.Range('D41').Style = $cleanStyle
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.21%  '
$ws.Range('E41').Style = $cleanStyle
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -3.02%  '
$ws.Range('E42').Style = $cleanStyle
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -3.86%  '
$ws.Range('E43').Style = $cleanStyle
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.65'
$ws.Range('D44').Style = $cleanStyle
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.63%  '
$ws.Range('E44').Style = $cleanStyle
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '103.96'
$ws.Range('D45').Style = $cleanStyle
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -2.74%  '
$ws.Range('E45').Style = $cleanStyle
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.97%  '
$ws.Range('E46').Style = $cleanStyle
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.37'
$ws.Range('D47').Style = $cleanStyle
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +3.00%  '
$ws.Range('E47').Style = $cleanStyle
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.11'
$ws.Range('D48').Style = $cleanStyle
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.62%  '
$ws.Range('E48').Style = $cleanStyle
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.91%  '
$ws.Range('E49').Style = $cleanStyle
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.69'
$ws.Range('D50').Style = $cleanStyle
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.19%  '
$ws.Range('E50').Style = $cleanStyle
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.420.94'
$ws.Range('D51').Style = $cleanStyle
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.55%  '
$ws.Range('E51').Style = $cleanStyle
